$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  46034, 13467.6826692334, 13064.8206321092, 22395.86, 9272.75768704457, -2.42840336859354),
    @(3,  46035, 13050.6969231705, 12433.6618819835, 12075.86, 8927.23855848847, 386.876685019665),
    @(4,  46036, 12578.754289796,  11882.0073088843, 12075.86, 8556.77742433483, 348.455197217465),
    @(5,  46037, 12274.9811196162, 11618.8009000478, 12075.86, 8330.97152608517, 328.079684422208),
    @(6,  46038, 11932.1866806542, 10679.4161981362, 12075.86, 8081.96366264703, 278.563327532635),
    @(7,  46039, 5247.16359875529, 7509.50426581102, 12075.86, 8468.05208663344, 162.570681351852),
    @(8,  46040, 5137.84212801753, 7298.53616501798, 12075.86, 8460.20612956141, 153.453428940808),
    @(9,  46041, 12719.4875548453, 11235.6512119228, 12075.86, 8648.2987964125,  325.337083680636),
    @(10, 46042, 12719.4875548453, 11452.7359612612, 12075.86, 8648.2987964125,  334.382281569737),
    @(11, 46043, 12719.4875548453, 11414.2064028501, 12075.86, 8648.2987964125,  332.776883302607),
    @(12, 46044, 12719.4875548453, 11480.4991576811, 12075.86, 8648.2987964125,  335.539081420565),
    @(13, 46045, 12719.4875548453, 10715.2918061095, 12075.86, 8648.2987964125,  303.655441771749),
    @(14, 46046, 5102.57643525512, 7375.12469436407, 12075.86, 8255.60568464833, 148.119599125517),
    @(15, 46047, 4995.11681641318, 7496.11122740767, 12075.86, 8247.86625619628, 152.838228483498)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
    $ws.Cells.Item($r, 6).Value2 = $row[6]
}
